# Append two more trial rows ("9" and "10") to the "Data" sheet and the
# matching summary rows to the "Summary" sheet, as produced by a further
# (still uncomplete) run of the ALL/SHAPESHIFTER experiment.

$wb = $excel.ActiveWorkbook

$chromosome9 = "[10, 7, 11, 16, 128, 0.7680441, 0.4905956, 0.58614576, 0.6788291, 0.3286342, 0.45604947, 0.63224334, 0.58431226, 0.06409139, 0.4486112, 0.5923541, 0.86245406, 0.5781349, 0.7811539, 0.7021895, 0.0062287697, 0.40543425, 0.5332268, 0.9989835, 0.724732, 0.18392125, 0.89512044, 0.83834064, 0.31168103, 0.41456914, 0.35160974, 0.3846827, 0.40361124, 0.51614547, 0.6810796, 0.45001304, 0.36917803, 0.6831876, 0.70975894, 0.467517, 0.7663588, 0.6001911, 0.56685805, 0.74872684, 0.76030254, 0.15015125, 0.5063697, 0.17825802, 0.45337522, 0.020946769, 0.41181844, 0.53246135]"

# Columns A..AD of the "Data" sheet, in order, for the two new trial rows.
# Both new trials (iteration 9 and 10) share the same chromosome / config
# values and only differ in A (iteration index), C (best fitness) and D
# (population average).
$dataRow12 = @(9.0, 1109.9333, 334.3911, 235.51036, 1109.9333, 158.4918, 30.0, $chromosome9, 37.966667, 16.806164, 549.86664, 27.60901, 291.83334, 97.64647, 0.033333335, 0.18257418, 45.0, 9.055386, 45.0, 9.055386, 28.733334, 10.599068, 0.0, 0.0, 45.0, 9.055386, 0.2, 0.48423418, 66.3, 27.327074)
$dataRow13 = @(10.0, 1109.9333, 518.5467, 208.61398, 1109.9333, 158.4918, 30.0, $chromosome9, 37.966667, 16.806164, 549.86664, 27.60901, 291.83334, 97.64647, 0.033333335, 0.18257418, 45.0, 9.055386, 45.0, 9.055386, 28.733334, 10.599068, 0.0, 0.0, 45.0, 9.055386, 0.2, 0.48423418, 66.3, 27.327074)

$ws = $wb.Worksheets.Item("Data")
for ($col = 1; $col -le $dataRow12.Length; $col++) {
    $ws.Cells.Item(12, $col).Value = $dataRow12[$col - 1]
}
for ($col = 1; $col -le $dataRow13.Length; $col++) {
    $ws.Cells.Item(13, $col).Value = $dataRow13[$col - 1]
}

# Matching rows (Iteration, Best fitness, Population average) on the
# "Summary" sheet.
$summaryRow11 = @(9.0, 1109.9333, 334.3911)
$summaryRow12 = @(10.0, 1109.9333, 518.5467)

$summary = $wb.Worksheets.Item("Summary")
for ($col = 1; $col -le $summaryRow11.Length; $col++) {
    $summary.Cells.Item(11, $col).Value = $summaryRow11[$col - 1]
}
for ($col = 1; $col -le $summaryRow12.Length; $col++) {
    $summary.Cells.Item(12, $col).Value = $summaryRow12[$col - 1]
}
